$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.579.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.775.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.36'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.776.04'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.19'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.32'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.408.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.771.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.565.60'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.36%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.25%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.51'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.49'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.741'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.70%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +11.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.13'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.48%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.24%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.71'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.109'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.51%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.28%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '447.93'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.06'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.97%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.15'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -8.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.833.45'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.62%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '138.97'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.04'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +10.04%  '
